$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns (by number) whose contents get re-shuffled across rows 2..27,
# following the permutation captured in $map (destRow -> sourceRow, using
# the *original* values before any writes happen).
$cols = @(8,9,11,12,13,14,15,16,17,18,19)  # H,I,K,L,M,N,O,P,Q,R,S

# destRow -> sourceRow mapping (based on the original/unshuffled data)
$map = @{
    2 = 6
    3 = 7
    4 = 8
    5 = 9
    6 = 10
    7 = 11
    8 = 15
    9 = 13
    10 = 25
    11 = 16
    12 = 27
    13 = 20
    14 = 17
    15 = 3
    16 = 18
    17 = 4
    18 = 12
    19 = 2
    20 = 22
    21 = 23
    22 = 19
    23 = 24
    24 = 26
    25 = 21
    26 = 5
    27 = 14
}

# Snapshot all source rows first, since several rows participate in swap
# cycles and would otherwise be overwritten before being read.
$snapshot = @{}
for ($r = 2; $r -le 27; $r++) {
    $rowVals = @{}
    foreach ($c in $cols) {
        $rowVals[$c] = $ws.Cells.Item($r, $c).Value2
    }
    $snapshot[$r] = $rowVals
}

# Write the shuffled values back, and set trial_total (F) = trial_block (E).
for ($r = 2; $r -le 27; $r++) {
    $srcRow = $map[$r]
    $srcVals = $snapshot[$srcRow]
    foreach ($c in $cols) {
        $ws.Cells.Item($r, $c).Value2 = $srcVals[$c]
    }
    $ws.Cells.Item($r, 6).Value2 = $ws.Cells.Item($r, 5).Value2
}
